# Prepare test data for the "ValidLogin" test case and expose it as a
# worksheet containing a UserName/Password header row plus one row of
# credentials (admin / manager).
#
# The original workbook has a single sheet named "Sheet1" (sheetId=1)
# containing two cells: A1="username", A2="admin".
#
# Target layout (sheet renamed to "ValidLogin"):
#   A1 = "UserName"   B1 = "Password"
#   A2 = "admin"      B2 = "manager"

$wb = $excel.ActiveWorkbook

# Remember the name of the sheet that currently exists so we can remove it
# once the replacement sheet is ready. Looking it up by name (rather than
# holding on to the worksheet object) keeps the reference correct even
# after the sheet collection changes.
$originalName = $wb.Worksheets.Item(1).Name
$original = $wb.Worksheets.Item($originalName)

# Copy the existing sheet right after itself. Doing a copy (rather than
# simply adding a blank sheet) preserves formatting/namespace declarations
# and also causes Excel to allocate a fresh, incremented sheetId for the
# new sheet.
$original.Copy($null, $original)

# The copy becomes sheet index 2 (right after the original at index 1).
$copyName = $wb.Worksheets.Item(2).Name
$new = $wb.Worksheets.Item($copyName)
$new.Name = "ValidLogin"

# Fill in the test data.
$new.Range("A1").Value = "UserName"
$new.Range("B1").Value = "Password"
$new.Range("A2").Value = "admin"
$new.Range("B2").Value = "manager"

# Remove the now-superseded original sheet.
$wb.Worksheets.Item($originalName).Delete()

# Re-select the surviving sheet (by name, since indices shifted after the
# delete) and update the view to match the authored state: cell B3
# selected and the zoom level bumped to 175%.
$ws = $wb.Worksheets.Item("ValidLogin")
$ws.Activate()
$ws.Range("B3").Select()
$excel.ActiveWindow.Zoom = 175
